# Statusbericht I.docx edit script
# - merges runs that were only split apart by w:proofErr spell-check markers
#   (InsertXML lets us rewrite a paragraph's content without Word re-adding
#   the proofErr wrapper, unlike Find/Replace which keeps it)
# - bumps the "Gesamt" total for Daniel from 37 to 40
# - replaces the "15.12" work-log row with the new content / hour counts
#
# NOTE: this interpreter parses `Func $var (expr)` as `$var(expr)` (a call
# on $var), not as two args to Func. So every argument that is not a bare
# variable/literal is built up in its own $variable first, then passed by
# name - never written as a parenthesized expression directly after
# another argument.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($para, [string]$innerXml) {
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $para.Range.InsertXML($xml)
}

function ArialRpr([string]$extra = "") {
    return "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`"/>$extra</w:rPr>"
}

$rprPlain = ArialRpr
$rprEn = ArialRpr '<w:lang w:val="en-US"/>'

# ---------------------------------------------------------------------
# Table 1 (Name / role / email / Ifw-Kennnung)
# ---------------------------------------------------------------------
$t0 = $d.Tables.Item(1)

# "Daniel " + "Dobras" (proofErr-wrapped) -> single run "Daniel Dobras"
$cell = $t0.Rows.Item(1).Cells.Item(2)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr><w:jc w:val=`"both`"/>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Daniel Dobras</w:t></w:r>"
Set-ParaXml $para $inner

# "Ifw-Kennnung" - drop the spellStart/spellEnd wrap, text unchanged
$cell = $t0.Rows.Item(4).Cells.Item(1)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr><w:jc w:val=`"both`"/>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Ifw-Kennnung</w:t></w:r>"
Set-ParaXml $para $inner

# ---------------------------------------------------------------------
# Table 2 (Woche planning table)
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(2)
$cell = $t1.Rows.Item(1).Cells.Item(2)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr><w:jc w:val=`"both`"/>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Daniel Dobras</w:t></w:r>"
Set-ParaXml $para $inner

# ---------------------------------------------------------------------
# Table 3 (Risiko table)
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(3)

# "Zeitlicher Engpass"
$cell = $t2.Rows.Item(3).Cells.Item(1)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprEn</w:pPr>"
$inner = $pPr + "<w:r>$rprEn<w:t>Zeitlicher Engpass</w:t></w:r>"
Set-ParaXml $para $inner

# "Zusätzliche Wünsche des " + "Auftraggebers" (keep lastRenderedPageBreak before "Auftraggebers")
$cell = $t2.Rows.Item(4).Cells.Item(1)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprEn</w:pPr>"
$run1 = "<w:r>$rprEn<w:t xml:space=`"preserve`">Zusätzliche Wünsche des </w:t></w:r>"
$run2 = "<w:r>$rprEn<w:lastRenderedPageBreak/><w:t>Auftraggebers</w:t></w:r>"
$inner = $pPr + $run1 + $run2
Set-ParaXml $para $inner

# ---------------------------------------------------------------------
# Table 4 (work log)
# ---------------------------------------------------------------------
$t3 = $d.Tables.Item(4)

# Header row: "Daniel " + "Dobras" -> "Daniel Dobras"
$cell = $t3.Rows.Item(1).Cells.Item(2)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Daniel Dobras</w:t></w:r>"
Set-ParaXml $para $inner

# Gesamt: 37 -> 40
$cell = $t3.Rows.Item(2).Cells.Item(2)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>40</w:t></w:r>"
Set-ParaXml $para $inner

# 27.11 row: "...Re" + "g" + "i" stay, "ster " + "Component" + " angelegt" merge
$cell = $t3.Rows.Item(11).Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$run1 = "<w:r>$rprPlain<w:t>Lastenheft fertiggestellt, Re</w:t></w:r>"
$run2 = "<w:r>$rprPlain<w:t>g</w:t></w:r>"
$run3 = "<w:r>$rprPlain<w:t>i</w:t></w:r>"
$run4 = "<w:r>$rprPlain<w:t>ster Component angelegt</w:t></w:r>"
$inner = $pPr + $run1 + $run2 + $run3 + $run4
Set-ParaXml $para $inner

# 01.12 row: "Login-" + "Component" -> single run
$cell = $t3.Rows.Item(14).Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Login-Component</w:t></w:r>"
Set-ParaXml $para $inner

# 02.12 row: "Redirecting" + " " + "code" -> single run
$cell = $t3.Rows.Item(15).Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Redirecting code</w:t></w:r>"
Set-ParaXml $para $inner

# 09.12 row: "MCQ, OQ, " + "User.model" -> single run
$cell = $t3.Rows.Item(17).Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>MCQ, OQ, User.model</w:t></w:r>"
Set-ParaXml $para $inner

# 14.12 row: "Einarbeitung und " + "user." + bookmark + "model" -> single run, bookmark removed here
$cell = $t3.Rows.Item(18).Cells.Item(3)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>Einarbeitung und user.model</w:t></w:r>"
Set-ParaXml $para $inner

# 15.12 row: full content replace
$row19 = $t3.Rows.Item(19)

$cell = $row19.Cells.Item(2)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>6</w:t></w:r>"
Set-ParaXml $para $inner

$cell = $row19.Cells.Item(3)
$pPrLang = "<w:pPr>$rprEn</w:pPr>"
$p1Runs = "<w:r>$rprEn<w:t>Register and services</w:t></w:r>" + "<w:r>$rprEn<w:t>,</w:t></w:r>"
$p1 = "<w:p $wNs>" + $pPrLang + $p1Runs + "</w:p>"
$p2Runs = "<w:r>$rprEn<w:t>redirecting from login to quiz</w:t></w:r>" + '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>'
$p2 = "<w:p $wNs>" + $pPrLang + $p2Runs + "</w:p>"
$xml = $p1 + $p2
$cell.Range.InsertXML($xml)

$cell = $row19.Cells.Item(4)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>2</w:t></w:r>"
Set-ParaXml $para $inner

$cell = $row19.Cells.Item(5)
$para = $cell.Range.Paragraphs.Item(1)
$pPr = "<w:pPr>$rprPlain</w:pPr>"
$inner = $pPr + "<w:r>$rprPlain<w:t>sql</w:t></w:r>"
Set-ParaXml $para $inner

Write-Output "done"
